$d = $word.ActiveDocument

# --- Body text: "A TERE," -> "A QWER," (bold run in document.xml) ---
$bodyRng = $d.Content
$bodyRng.Find.Execute("TERE", $true, $false, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# --- Header (primary) text replacements, in document order ---
$sec = $d.Sections.First
$hdr = $sec.Headers.Item(1)
$hdrRng = $hdr.Range

# "DIRETORIA DE ENSINO REGIAO TRE" -> "... QWER"
$hdrRng.Find.Execute("TRE", $true, $false, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# "TERE - DEP." -> "QWER - DEP."
$hdrRng.Find.Execute("TERE", $true, $false, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# Address line: five "Tre" occurrences -> Qwer, Qwer, Qewr, Qewr, Qwer
$addressReplacements = @("Qwer", "Qwer", "Qewr", "Qewr", "Qwer")
foreach ($rep in $addressReplacements) {
    $hdrRng.Find.Execute("Tre", $true, $false, $false, $false, $false, $true, 1, $false, $rep, 1) | Out-Null
}

# CEP / Tel line: two "tre" occurrences -> qwer, qwer
$hdrRng.Find.Execute("tre", $true, $false, $false, $false, $false, $true, 1, $false, "qwer", 1) | Out-Null
$hdrRng.Find.Execute("tre", $true, $false, $false, $false, $false, $true, 1, $false, "qwer", 1) | Out-Null

# Email line: "tre" -> "qwer"
$hdrRng.Find.Execute("tre", $true, $false, $false, $false, $false, $true, 1, $false, "qwer", 1) | Out-Null
